$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the two column header labels that changed as part of the
# Meerkat DB changes (Person_ID -> PersonBusinessKey, Role_ID -> RoleBusinessKey)
$ws.Range("C2").Value = "PersonBusinessKey"
$ws.Range("D2").Value = "RoleBusinessKey"
